# ------------------------------------------------------------------
# Edit script: update species names, add new Asclepias incarnata rows
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RFriendly")

# --- Update species labels to their full scientific names ---
# (the underlying shared string for "Dalea" is retired and a new
#  "Dalea candida" string is introduced; likewise for the others)
$ws.Range("A2:A19").Value   = "Lespedeza capitata"
$ws.Range("A20:A27").Value  = "Dalea candida"
$ws.Range("A28:A39").Value  = "Asclepias tuberosa"
$ws.Range("A40:A127").Value = "Amorpha canescens"
$ws.Range("A128:A214").Value= "Baptisia bracteata"

# --- Append new data rows 215-300 for a new species: Asclepias incarnata ---
$newAB = New-Object "object[,]" 86,2
$newGH = New-Object "object[,]" 86,2
$newAB[0,0] = "Asclepias incarnata"; $newAB[0,1] = "Low"; $newGH[0,0] = 0.6; $newGH[0,1] = 0.5
$newAB[1,0] = "Asclepias incarnata"; $newAB[1,1] = "Low"; $newGH[1,0] = 0.255; $newGH[1,1] = 0.32
$newAB[2,0] = "Asclepias incarnata"; $newAB[2,1] = "Low"; $newGH[2,0] = 0.545; $newGH[2,1] = 0.775
$newAB[3,0] = "Asclepias incarnata"; $newAB[3,1] = "Low"; $newGH[3,0] = 0.305; $newGH[3,1] = 0.24
$newAB[4,0] = "Asclepias incarnata"; $newAB[4,1] = "Low"; $newGH[4,0] = 0.56; $newGH[4,1] = 0.535
$newAB[5,0] = "Asclepias incarnata"; $newAB[5,1] = "Low"; $newGH[5,0] = 0.36; $newGH[5,1] = 0.425
$newAB[6,0] = "Asclepias incarnata"; $newAB[6,1] = "Low"; $newGH[6,0] = 0.485; $newGH[6,1] = 0.385
$newAB[7,0] = "Asclepias incarnata"; $newAB[7,1] = "Low"; $newGH[7,0] = 0.5; $newGH[7,1] = 0.45
$newAB[8,0] = "Asclepias incarnata"; $newAB[8,1] = "Low"; $newGH[8,0] = 0.365; $newGH[8,1] = 0.475
$newAB[9,0] = "Asclepias incarnata"; $newAB[9,1] = "Low"; $newGH[9,0] = 0.325; $newGH[9,1] = 0.375
$newAB[10,0] = "Asclepias incarnata"; $newAB[10,1] = "Low"; $newGH[10,0] = 0.455; $newGH[10,1] = 0.4
$newAB[11,0] = "Asclepias incarnata"; $newAB[11,1] = "Low"; $newGH[11,0] = 0.595; $newGH[11,1] = 0.565
$newAB[12,0] = "Asclepias incarnata"; $newAB[12,1] = "Low"; $newGH[12,0] = 0.265; $newGH[12,1] = 0.29
$newAB[13,0] = "Asclepias incarnata"; $newAB[13,1] = "Low"; $newGH[13,0] = 0.4; $newGH[13,1] = 0.305
$newAB[14,0] = "Asclepias incarnata"; $newAB[14,1] = "Low"; $newGH[14,0] = 0.45; $newGH[14,1] = 0.365
$newAB[15,0] = "Asclepias incarnata"; $newAB[15,1] = "Low"; $newGH[15,0] = 0.575; $newGH[15,1] = 0.355
$newAB[16,0] = "Asclepias incarnata"; $newAB[16,1] = "Low"; $newGH[16,0] = 0.645; $newGH[16,1] = 0.5
$newAB[17,0] = "Asclepias incarnata"; $newAB[17,1] = "Low"; $newGH[17,0] = 0.43; $newGH[17,1] = 0.265
$newAB[18,0] = "Asclepias incarnata"; $newAB[18,1] = "Low"; $newGH[18,0] = 0.28; $newGH[18,1] = 0.365
$newAB[19,0] = "Asclepias incarnata"; $newAB[19,1] = "Low"; $newGH[19,0] = 0.35; $newGH[19,1] = 0.365
$newAB[20,0] = "Asclepias incarnata"; $newAB[20,1] = "Low"; $newGH[20,0] = 0.32; $newGH[20,1] = 0.25
$newAB[21,0] = "Asclepias incarnata"; $newAB[21,1] = "Low"; $newGH[21,0] = 0.645; $newGH[21,1] = 0.515
$newAB[22,0] = "Asclepias incarnata"; $newAB[22,1] = "Low"; $newGH[22,0] = 0.445; $newGH[22,1] = 0.575
$newAB[23,0] = "Asclepias incarnata"; $newAB[23,1] = "Low"; $newGH[23,0] = 0.4; $newGH[23,1] = 0.33
$newAB[24,0] = "Asclepias incarnata"; $newAB[24,1] = "Low"; $newGH[24,0] = 0.72; $newGH[24,1] = 0.735
$newAB[25,0] = "Asclepias incarnata"; $newAB[25,1] = "Low"; $newGH[25,0] = 0.555; $newGH[25,1] = 0.52
$newAB[26,0] = "Asclepias incarnata"; $newAB[26,1] = "Low"; $newGH[26,0] = 0.59; $newGH[26,1] = 0.535
$newAB[27,0] = "Asclepias incarnata"; $newAB[27,1] = "Low"; $newGH[27,0] = 0.36; $newGH[27,1] = 0.365
$newAB[28,0] = "Asclepias incarnata"; $newAB[28,1] = "Low"; $newGH[28,0] = 0.775; $newGH[28,1] = 0.54
$newAB[29,0] = "Asclepias incarnata"; $newAB[29,1] = "Low"; $newGH[29,0] = 0.35; $newGH[29,1] = 0.555
$newAB[30,0] = "Asclepias incarnata"; $newAB[30,1] = "Low"; $newGH[30,0] = 0.34; $newGH[30,1] = 0.43
$newAB[31,0] = "Asclepias incarnata"; $newAB[31,1] = "Low"; $newGH[31,0] = 0.31; $newGH[31,1] = 0.5
$newAB[32,0] = "Asclepias incarnata"; $newAB[32,1] = "Low"; $newGH[32,0] = 0.14; $newGH[32,1] = 0.335
$newAB[33,0] = "Asclepias incarnata"; $newAB[33,1] = "Low"; $newGH[33,0] = 0.64; $newGH[33,1] = 0.915
$newAB[34,0] = "Asclepias incarnata"; $newAB[34,1] = "Low"; $newGH[34,0] = 0.465; $newGH[34,1] = 0.515
$newAB[35,0] = "Asclepias incarnata"; $newAB[35,1] = "Low"; $newGH[35,0] = 0.545; $newGH[35,1] = 0.565
$newAB[36,0] = "Asclepias incarnata"; $newAB[36,1] = "Low"; $newGH[36,0] = 0.56; $newGH[36,1] = 0.62
$newAB[37,0] = "Asclepias incarnata"; $newAB[37,1] = "Low"; $newGH[37,0] = 0.445; $newGH[37,1] = 0.42
$newAB[38,0] = "Asclepias incarnata"; $newAB[38,1] = "Low"; $newGH[38,0] = 0.5; $newGH[38,1] = 1.08
$newAB[39,0] = "Asclepias incarnata"; $newAB[39,1] = "Low"; $newGH[39,0] = 0.275; $newGH[39,1] = 0.41
$newAB[40,0] = "Asclepias incarnata"; $newAB[40,1] = "Low"; $newGH[40,0] = 0.43; $newGH[40,1] = 0.665
$newAB[41,0] = "Asclepias incarnata"; $newAB[41,1] = "Low"; $newGH[41,0] = 0.45; $newGH[41,1] = 0.405
$newAB[42,0] = "Asclepias incarnata"; $newAB[42,1] = "Low"; $newGH[42,0] = 0.43; $newGH[42,1] = 0.345
$newAB[43,0] = "Asclepias incarnata"; $newAB[43,1] = "Low"; $newGH[43,0] = 0.63; $newGH[43,1] = 0.41
$newAB[44,0] = "Asclepias incarnata"; $newAB[44,1] = "Low"; $newGH[44,0] = 0.415; $newGH[44,1] = 0.475
$newAB[45,0] = "Asclepias incarnata"; $newAB[45,1] = "Low"; $newGH[45,0] = 0.15; $newGH[45,1] = 0.33
$newAB[46,0] = "Asclepias incarnata"; $newAB[46,1] = "Low"; $newGH[46,0] = 0.48; $newGH[46,1] = 0.525
$newAB[47,0] = "Asclepias incarnata"; $newAB[47,1] = "Low"; $newGH[47,0] = 0.375; $newGH[47,1] = 0.355
$newAB[48,0] = "Asclepias incarnata"; $newAB[48,1] = "Low"; $newGH[48,0] = 0.395; $newGH[48,1] = 0.455
$newAB[49,0] = "Asclepias incarnata"; $newAB[49,1] = "High"; $newGH[49,0] = 0.545; $newGH[49,1] = 0.445
$newAB[50,0] = "Asclepias incarnata"; $newAB[50,1] = "High"; $newGH[50,0] = 0.5; $newGH[50,1] = 0.845
$newAB[51,0] = "Asclepias incarnata"; $newAB[51,1] = "High"; $newGH[51,0] = 0.295; $newGH[51,1] = 0.355
$newAB[52,0] = "Asclepias incarnata"; $newAB[52,1] = "High"; $newGH[52,0] = 0.55; $newGH[52,1] = 0.88
$newAB[53,0] = "Asclepias incarnata"; $newAB[53,1] = "High"; $newGH[53,0] = 0.475; $newGH[53,1] = 0.625
$newAB[54,0] = "Asclepias incarnata"; $newAB[54,1] = "High"; $newGH[54,0] = 0.225; $newGH[54,1] = 0.555
$newAB[55,0] = "Asclepias incarnata"; $newAB[55,1] = "High"; $newGH[55,0] = 0.15; $newGH[55,1] = 0.32
$newAB[56,0] = "Asclepias incarnata"; $newAB[56,1] = "High"; $newGH[56,0] = 0.61; $newGH[56,1] = 0.515
$newAB[57,0] = "Asclepias incarnata"; $newAB[57,1] = "High"; $newGH[57,0] = 0.54; $newGH[57,1] = 0.45
$newAB[58,0] = "Asclepias incarnata"; $newAB[58,1] = "High"; $newGH[58,0] = 0.775; $newGH[58,1] = 0.82
$newAB[59,0] = "Asclepias incarnata"; $newAB[59,1] = "High"; $newGH[59,0] = 0.455; $newGH[59,1] = 0.86
$newAB[60,0] = "Asclepias incarnata"; $newAB[60,1] = "High"; $newGH[60,0] = 0.135; $newGH[60,1] = 0.415
$newAB[61,0] = "Asclepias incarnata"; $newAB[61,1] = "High"; $newGH[61,0] = 0.25; $newGH[61,1] = 0.465
$newAB[62,0] = "Asclepias incarnata"; $newAB[62,1] = "High"; $newGH[62,0] = 0.57; $newGH[62,1] = 0.56
$newAB[63,0] = "Asclepias incarnata"; $newAB[63,1] = "High"; $newGH[63,0] = 0.375; $newGH[63,1] = 0.55
$newAB[64,0] = "Asclepias incarnata"; $newAB[64,1] = "High"; $newGH[64,0] = 0.395; $newGH[64,1] = 0.68
$newAB[65,0] = "Asclepias incarnata"; $newAB[65,1] = "High"; $newGH[65,0] = 0.26; $newGH[65,1] = 0.545
$newAB[66,0] = "Asclepias incarnata"; $newAB[66,1] = "High"; $newGH[66,0] = 0.635; $newGH[66,1] = 0.5
$newAB[67,0] = "Asclepias incarnata"; $newAB[67,1] = "High"; $newGH[67,0] = 0.55; $newGH[67,1] = 0.625
$newAB[68,0] = "Asclepias incarnata"; $newAB[68,1] = "High"; $newGH[68,0] = 0.04; $newGH[68,1] = 0.175
$newAB[69,0] = "Asclepias incarnata"; $newAB[69,1] = "High"; $newGH[69,0] = 0.66; $newGH[69,1] = 0.555
$newAB[70,0] = "Asclepias incarnata"; $newAB[70,1] = "High"; $newGH[70,0] = 0.395; $newGH[70,1] = 0.665
$newAB[71,0] = "Asclepias incarnata"; $newAB[71,1] = "High"; $newGH[71,0] = 0.725; $newGH[71,1] = 0.84
$newAB[72,0] = "Asclepias incarnata"; $newAB[72,1] = "High"; $newGH[72,0] = 0.5; $newGH[72,1] = 0.625
$newAB[73,0] = "Asclepias incarnata"; $newAB[73,1] = "High"; $newGH[73,0] = 0.78; $newGH[73,1] = 0.67
$newAB[74,0] = "Asclepias incarnata"; $newAB[74,1] = "High"; $newGH[74,0] = 0.255; $newGH[74,1] = 0.475
$newAB[75,0] = "Asclepias incarnata"; $newAB[75,1] = "High"; $newGH[75,0] = 0.515; $newGH[75,1] = 0.535
$newAB[76,0] = "Asclepias incarnata"; $newAB[76,1] = "High"; $newGH[76,0] = 0.62; $newGH[76,1] = 1.065
$newAB[77,0] = "Asclepias incarnata"; $newAB[77,1] = "High"; $newGH[77,0] = 0.13; $newGH[77,1] = 0.285
$newAB[78,0] = "Asclepias incarnata"; $newAB[78,1] = "High"; $newGH[78,0] = 0.675; $newGH[78,1] = 0.91
$newAB[79,0] = "Asclepias incarnata"; $newAB[79,1] = "High"; $newGH[79,0] = 0.67; $newGH[79,1] = 0.725
$newAB[80,0] = "Asclepias incarnata"; $newAB[80,1] = "High"; $newGH[80,0] = 0.335; $newGH[80,1] = 0.445
$newAB[81,0] = "Asclepias incarnata"; $newAB[81,1] = "High"; $newGH[81,0] = 0.64; $newGH[81,1] = 0.595
$newAB[82,0] = "Asclepias incarnata"; $newAB[82,1] = "High"; $newGH[82,0] = 0.5; $newGH[82,1] = 0.485
$newAB[83,0] = "Asclepias incarnata"; $newAB[83,1] = "High"; $newGH[83,0] = 0.305; $newGH[83,1] = 0.595
$newAB[84,0] = "Asclepias incarnata"; $newAB[84,1] = "High"; $newGH[84,0] = 0.355; $newGH[84,1] = 0.66
$newAB[85,0] = "Asclepias incarnata"; $newAB[85,1] = "High"; $newGH[85,0] = 0.355; $newGH[85,1] = 0.445

$ws.Range("A215:B300").Value = $newAB
$ws.Range("G215:H300").Value = $newGH

# --- Update the view state: scroll position & active selection ---
$ws.Activate()
$ws.Range("B284").Select()
